# Apply commit: insert a new headline ("VPBS tiep tuc giai chap...") into
# the list and append a batch of further headlines, per the author's diff.
#
# Final A1:A17 headline list (after the edit):
#   A1  Tổng thống Mỹ ký ban hành luật TPA, mở đường cho TPP
#   A2  Quảng Ninh "nhượng bộ", đền bù đợt 1 cho CTCP Quốc tế Hoàng Gia (RIC) 110 tỷ đồng
#   A3  VPBS tiếp tục giải chấp thêm cổ phiếu JVC của ông Lê Văn Hướng
#   A4  Quảng Ninh "nhượng bộ", đền bù đợt 1 cho CTCP Quốc tế Hoàng Gia (RIC) 110 tỷ đồng
#   A5  Công ty chứng khoán được mở room 100% ngay
#   A6  Quy định "cong vênh" khiến xuất khẩu mất hàng tỷ đô
#   A7  TÔI ĐẦU TƯ: Muốn thắng thị trường, trước tiên tuyệt đối chỉ mua cổ phiếu tại giá trần
#   A8  'Giấc mộng châu Âu' có tan vỡ vì Hy Lạp?
#   A9  Từ hôm nay (1/7) Luật nhà ở, Luật kinh doanh BĐS 2014 chính thức có hiệu lực
#   A10 Tổng thống Mỹ ký ban hành luật TPA, mở đường cho TPP
#   A11 Nhịp đập Thị trường 01/07: Thanh khoản sụt giảm mạnh
#   A12 JVC: Ông Lê Văn Hướng tiếp tục bị bán giải chấp
#   A13 Vì sao HHS giảm sàn?
#   A14 Những khoản cổ tức chốt quyền trong nửa đầu tháng 7
#   A15 PMI sản xuất tháng 6 giảm nhưng vẫn duy trì trên mốc 50 điểm
#   A16 Nhịp đập Thị trường 01/07: Tiền chuyển hướng qua một số cổ phiếu nhỏ
#   A17 Hy Lạp vỡ nợ, EU từ chối cứu trợ

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$headlines = @(
    "Tổng thống Mỹ ký ban hành luật TPA, mở đường cho TPP",
    "Quảng Ninh `"nhượng bộ`", đền bù đợt 1 cho CTCP Quốc tế Hoàng Gia (RIC) 110 tỷ đồng",
    "VPBS tiếp tục giải chấp thêm cổ phiếu JVC của ông Lê Văn Hướng",
    "Quảng Ninh `"nhượng bộ`", đền bù đợt 1 cho CTCP Quốc tế Hoàng Gia (RIC) 110 tỷ đồng",
    "Công ty chứng khoán được mở room 100% ngay",
    "Quy định “cong vênh” khiến xuất khẩu mất hàng tỷ đô",
    "TÔI ĐẦU TƯ: Muốn thắng thị trường, trước tiên tuyệt đối chỉ mua cổ phiếu tại giá trần",
    "''Giấc mộng châu Âu' có tan vỡ vì Hy Lạp?",
    "Từ hôm nay (1/7) Luật nhà ở, Luật kinh doanh BĐS 2014 chính thức có hiệu lực",
    "Tổng thống Mỹ ký ban hành luật TPA, mở đường cho TPP",
    "Nhịp đập Thị trường 01/07: Thanh khoản sụt giảm mạnh",
    "JVC: Ông Lê Văn Hướng tiếp tục bị bán giải chấp",
    "Vì sao HHS giảm sàn?",
    "Những khoản cổ tức chốt quyền trong nửa đầu tháng 7",
    "PMI sản xuất tháng 6 giảm nhưng vẫn duy trì trên mốc 50 điểm",
    "Nhịp đập Thị trường 01/07: Tiền chuyển hướng qua một số cổ phiếu nhỏ",
    "Hy Lạp vỡ nợ, EU từ chối cứu trợ"
)

for ($i = 0; $i -lt $headlines.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $headlines[$i]
}

# The headline at A8 starts with a literal apostrophe; Excel treats a
# leading "'" in an assigned .Value as a force-text prefix marker (stripped
# from the stored string, cell flagged quotePrefix). Doubling it above
# preserved the real apostrophe in the text; resetting the style here drops
# the extraneous quotePrefix formatting Excel applied along the way.
$ws.Range("A8").Style = "Normal"
